$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph right before
# the final section properties. Insert a brand new paragraph right after
# it (i.e. still before the sectPr) to hold the new commentary text.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Alignment = 3   # wdAlignParagraphJustify -> <w:jc w:val="both"/>

$sentenceOne = "Природа графиков. Экстремальные точки связаны с тем, что при потоках большем 8 накладные расходы на «сворачивание» потоков занимают больше времени чем экономит времени наше ускорение за счет большого количества потоков."
$sentenceTwo = " Экстремальная точка в восьми потоках связана с тем, что при этом числе достигается наибольшее ускорение работы алгоритма."

$newRange = $newPara.Range
$newRange.InsertAfter($sentenceOne + $sentenceTwo)

# The paragraph currently holds the whole sentence in a single run. Touch
# the formatting of just the first sentence (on, then back off) so Word
# keeps it as its own run, matching the two separate <w:r> runs used by
# the author (both runs end up with identical, unmodified formatting).
$paraStart = $newPara.Range.Start
$firstRunRange = $d.Range($paraStart, $paraStart + $sentenceOne.Length)
$firstRunRange.Font.Bold = 1
$firstRunRange.Font.Bold = 0

Write-Output "Inserted commentary paragraph."
